$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8-11 down to 9-12.
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new weekly price entry.
$ws.Range("A8").Value = 11
$ws.Range("B8").Value = "Vega Monumental Concepción"
$ws.Range("C8").Value = "Bíobío"
$ws.Range("D8").Value = 44839
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 300000000
$ws.Range("G8").Value = "Espárragos"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 1700
$ws.Range("L8").Value = 1800
$ws.Range("M8").Value = 1760
$ws.Range("N8").Value = "$/kilo"
$ws.Range("O8").Value = "Provincia de Linares"
$ws.Range("P8").Value = 1760
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = "Hortaliza"
